$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------------
# This script applies the weekly crime-data refresh described by the commit:
#   - Bumps the report volume/week number and the covered date range (header).
#   - Refreshes the statistics table (rows 14-29) with the newly collected
#     weekly / 28-day / YTD / 2-year figures and their computed %% changes.
# Some cells flip between a numeric value and a text placeholder ("0" or
# "***.*") used by this report for suppressed/undefined values; for those we
# copy formatting from an untouched donor cell that already has the desired
# style before writing the new value, so no new cell styles are created.
# ----------------------------------------------------------------------------

# --- Cells that switch between numeric and text-placeholder styles ---
# (copy style/format from a stable, unmodified donor cell first)
$ws.Range("C30").Copy($ws.Range("F14"))
$ws.Range("C30").Copy($ws.Range("D15"))
$ws.Range("E30").Copy($ws.Range("E15"))
$ws.Range("C30").Copy($ws.Range("F15"))
$ws.Range("C30").Copy($ws.Range("C22"))
$ws.Range("C30").Copy($ws.Range("D26"))
$ws.Range("E30").Copy($ws.Range("E26"))
$ws.Range("J30").Copy($ws.Range("C27"))
$ws.Range("C30").Copy($ws.Range("D27"))
$ws.Range("E30").Copy($ws.Range("E27"))
$ws.Range("J30").Copy($ws.Range("C28"))
$ws.Range("C30").Copy($ws.Range("D28"))
$ws.Range("E30").Copy($ws.Range("E28"))
$ws.Range("J30").Copy($ws.Range("F28"))
$ws.Range("J30").Copy($ws.Range("C29"))
$ws.Range("C30").Copy($ws.Range("D29"))
$ws.Range("E30").Copy($ws.Range("E29"))
$ws.Range("J30").Copy($ws.Range("F29"))

# --- Set the new values for those same cells ---
$ws.Range("F14").Value = "0"
$ws.Range("D15").Value = "0"
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = "0"
$ws.Range("C22").Value = "0"
$ws.Range("D26").Value = "0"
$ws.Range("E26").Value = "***.*"
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = "0"
$ws.Range("E27").Value = "***.*"
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = "0"
$ws.Range("E28").Value = "***.*"
$ws.Range("F28").Value = 1
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = "0"
$ws.Range("E29").Value = "***.*"
$ws.Range("F29").Value = 1

# --- Plain value updates (style/number format unchanged) ---
$ws.Range("H15").Value = -100
$ws.Range("L15").Value = 37.931034482758
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = 62.5
$ws.Range("F16").Value = 51
$ws.Range("G16").Value = 44
$ws.Range("H16").Value = 15.909090909090
$ws.Range("I16").Value = 666
$ws.Range("J16").Value = 558
$ws.Range("K16").Value = 19.354838709677
$ws.Range("L16").Value = 52.402745995423
$ws.Range("M16").Value = 49.327354260089
$ws.Range("N16").Value = -65.094339622641
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = 22
$ws.Range("E17").Value = -27.272727272727
$ws.Range("F17").Value = 80
$ws.Range("G17").Value = 68
$ws.Range("H17").Value = 17.647058823529
$ws.Range("I17").Value = 1024
$ws.Range("J17").Value = 853
$ws.Range("K17").Value = 20.046893317702
$ws.Range("L17").Value = 44.022503516174
$ws.Range("M17").Value = 134.324942791762
$ws.Range("N17").Value = -5.360443622920
$ws.Range("C18").Value = 14
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 250
$ws.Range("F18").Value = 30
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = 30.434782608695
$ws.Range("I18").Value = 327
$ws.Range("J18").Value = 365
$ws.Range("K18").Value = -10.410958904109
$ws.Range("L18").Value = 65.989847715736
$ws.Range("M18").Value = 62.686567164179
$ws.Range("N18").Value = -76.026392961876
$ws.Range("C19").Value = 16
$ws.Range("E19").Value = -5.882352941176
$ws.Range("F19").Value = 62
$ws.Range("G19").Value = 61
$ws.Range("H19").Value = 1.639344262295
$ws.Range("I19").Value = 749
$ws.Range("J19").Value = 755
$ws.Range("K19").Value = -0.794701986754
$ws.Range("L19").Value = 1.216216216216
$ws.Range("M19").Value = 72.183908045977
$ws.Range("N19").Value = -2.219321148825
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -14.285714285714
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 44.444444444444
$ws.Range("I20").Value = 340
$ws.Range("J20").Value = 305
$ws.Range("K20").Value = 11.475409836065
$ws.Range("L20").Value = 69.154228855721
$ws.Range("M20").Value = 169.84126984127
$ws.Range("N20").Value = -47.932618683001
$ws.Range("C21").Value = 65
$ws.Range("D21").Value = 58
$ws.Range("E21").Value = 12.068965517241
$ws.Range("F21").Value = 249
$ws.Range("G21").Value = 217
$ws.Range("H21").Value = 14.746543778801
$ws.Range("I21").Value = 3158
$ws.Range("J21").Value = 2882
$ws.Range("K21").Value = 9.576682859125
$ws.Range("L21").Value = 35.420240137221
$ws.Range("M21").Value = 87.752675386444
$ws.Range("N21").Value = -46.682424447070
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -20
$ws.Range("J22").Value = 85
$ws.Range("K22").Value = -21.176470588235
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 11
$ws.Range("E23").Value = -18.181818181818
$ws.Range("G23").Value = 37
$ws.Range("H23").Value = 2.702702702702
$ws.Range("I23").Value = 488
$ws.Range("J23").Value = 402
$ws.Range("K23").Value = 21.393034825870
$ws.Range("L23").Value = 56.913183279742
$ws.Range("M23").Value = 65.986394557823
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 51
$ws.Range("E24").Value = -43.137254901960
$ws.Range("F24").Value = 150
$ws.Range("G24").Value = 169
$ws.Range("H24").Value = -11.242603550295
$ws.Range("I24").Value = 1642
$ws.Range("J24").Value = 1952
$ws.Range("K24").Value = -15.881147540983
$ws.Range("L24").Value = 15.796897038081
$ws.Range("M24").Value = 15.147265077138
$ws.Range("C25").Value = 33
$ws.Range("D25").Value = 33
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 92
$ws.Range("G25").Value = 82
$ws.Range("H25").Value = 12.195121951219
$ws.Range("I25").Value = 1144
$ws.Range("J25").Value = 1055
$ws.Range("K25").Value = 8.436018957345
$ws.Range("L25").Value = 20.802534318901
$ws.Range("M25").Value = 0.615655233069
$ws.Range("F26").Value = 2
$ws.Range("H26").Value = -66.666666666666
$ws.Range("L26").Value = -11.764705882352
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 106
$ws.Range("K27").Value = 21.839080459770
$ws.Range("L27").Value = 13.978494623655
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -87.5
$ws.Range("I28").Value = 35
$ws.Range("K28").Value = -46.969696969697
$ws.Range("L28").Value = -50.704225352112
$ws.Range("M28").Value = -43.548387096774
$ws.Range("N28").Value = -83.644859813084
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -75
$ws.Range("I29").Value = 33
$ws.Range("K29").Value = -36.538461538461
$ws.Range("L29").Value = -47.619047619047
$ws.Range("M29").Value = -37.735849056603
$ws.Range("N29").Value = -82.901554404145

# --- Header: bump report Volume/Number and the covered week date range ---
$ws.Range("A8").Value = "Volume 30   Number  52"
$ws.Range("C9").Value = "Report Covering the Week  12/25/2023  Through  12/31/2023"
